# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.344.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "'1.830.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "'314.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D8").Value = "'0.3700"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("D9").Value = "'0.07271"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").Value = "'0.8674"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("D12").Value = "'1.827.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "'6.743"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "'0.07098"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "'5.325"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").Value = "'89.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "'0.000008881"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.83%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'15.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.51%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "'27.364.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.142"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'10.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "'2.053.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.000"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'152.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.188"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.28%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'18.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'5.260"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.33%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'116.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.96%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.08881"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").Value = "'1.204"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.99%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7615"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.481"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.824"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "'1.006"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.124"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01985"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05283"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.275"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.886"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1703"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5091"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'8.721"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'10.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'108.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4782"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06399"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.674"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.844"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.79%  "
